$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2,1).Value = [double]"6.217406212239571e-21"
$ws.Cells.Item(2,2).Value = [double]"-5.275625057681557e-21"
$ws.Cells.Item(2,3).Value = [double]"3.219057089704902e-11"
$ws.Cells.Item(2,4).Value = [double]"5.275620846622326e-21"
$ws.Cells.Item(2,5).Value = [double]"2.731461193349775e-11"
$ws.Cells.Item(2,6).Value = [double]"-2.731459451213997e-11"
$ws.Cells.Item(2,7).Value = [double]"1.931432204919333e-10"
$ws.Cells.Item(3,1).Value = [double]"1.931437842995357e-10"
$ws.Cells.Item(3,2).Value = [double]"3.45647655392911e-16"
$ws.Cells.Item(3,3).Value = [double]"1"
$ws.Cells.Item(3,4).Value = [double]"-4.901593125158773e-16"
$ws.Cells.Item(3,5).Value = [double]"-1.365726115649133e-10"
$ws.Cells.Item(3,6).Value = [double]"-1.365734379330293e-10"
$ws.Cells.Item(3,7).Value = [double]"5.334071672239157e-17"
$ws.Cells.Item(3,8).Value = [double]"3.219057089704904e-11"
$ws.Cells.Item(4,1).Value = [double]"3.219059206605705e-11"
$ws.Cells.Item(4,2).Value = [double]"-2.731460402022783e-11"
$ws.Cells.Item(4,3).Value = [double]"2.459222016602375e-17"
$ws.Cells.Item(4,4).Value = [double]"2.731460631153943e-11"
$ws.Cells.Item(4,5).Value = [double]"1.030076663856088e-16"
$ws.Cells.Item(4,6).Value = [double]"-9.760978308568577e-17"
$ws.Cells.Item(4,7).Value = [double]"1"
$ws.Cells.Item(4,8).Value = [double]"1.931435135779936e-10"
$ws.Cells.Item(5,1).Value = [double]"1"
$ws.Cells.Item(5,2).Value = [double]"1.365730315576971e-10"
$ws.Cells.Item(5,3).Value = [double]"-1.931437241363401e-10"
$ws.Cells.Item(5,4).Value = [double]"1.365730201011393e-10"
$ws.Cells.Item(5,5).Value = [double]"2.156604161989864e-16"
$ws.Cells.Item(5,6).Value = [double]"2.775036610834975e-17"
$ws.Cells.Item(5,7).Value = [double]"-3.219057849715358e-11"
$ws.Cells.Item(5,8).Value = [double]"-6.217400698932025e-21"
$ws.Cells.Item(6,1).Value = [double]"-4.464578364049936e-17"
$ws.Cells.Item(6,2).Value = [double]"-1.365730201327928e-10"
$ws.Cells.Item(6,3).Value = [double]"-4.432244306581715e-16"
$ws.Cells.Item(6,4).Value = [double]"1.365730315260431e-10"
$ws.Cells.Item(6,5).Value = [double]"0.70710695654703"
$ws.Cells.Item(6,6).Value = [double]"-0.7071065052238935"
$ws.Cells.Item(6,7).Value = [double]"2.76967488568725e-17"
$ws.Cells.Item(6,8).Value = [double]"3.862868507645885e-11"
$ws.Cells.Item(7,1).Value = [double]"-6.806329781980839e-17"
$ws.Cells.Item(7,2).Value = [double]"-1.365730315260432e-10"
$ws.Cells.Item(7,3).Value = [double]"1.931434254926594e-10"
$ws.Cells.Item(7,4).Value = [double]"-1.365730201327929e-10"
$ws.Cells.Item(7,5).Value = [double]"0.7071066058260215"
$ws.Cells.Item(7,6).Value = [double]"0.7071070571490938"
$ws.Cells.Item(7,7).Value = [double]"-6.897807984476963e-21"
$ws.Cells.Item(7,8).Value = [double]"-2.735654737060307e-21"
$ws.Cells.Item(8,1).Value = [double]"4.352180292143494e-20"
$ws.Cells.Item(8,2).Value = [double]"-0.7071067516922724"
$ws.Cells.Item(8,3).Value = [double]"2.317721104587532e-10"
$ws.Cells.Item(8,4).Value = [double]"0.7071068106808213"
$ws.Cells.Item(8,5).Value = [double]"-1.365730596991424e-10"
$ws.Cells.Item(8,6).Value = [double]"1.365729725290458e-10"
$ws.Cells.Item(8,7).Value = [double]"-3.862868507645885e-11"
$ws.Cells.Item(8,8).Value = [double]"-7.460879960565068e-21"
$ws.Cells.Item(9,1).Value = [double]"1.931434253822942e-10"
$ws.Cells.Item(9,2).Value = [double]"-0.7071068106808214"
$ws.Cells.Item(9,3).Value = [double]"-2.984356843292659e-20"
$ws.Cells.Item(9,4).Value = [double]"-0.7071067516922726"
$ws.Cells.Item(9,5).Value = [double]"-1.365729919992527e-10"
$ws.Cells.Item(9,6).Value = [double]"-1.365730790902023e-10"
$ws.Cells.Item(9,7).Value = [double]"4.973915112249691e-21"
$ws.Cells.Item(9,8).Value = [double]"7.685437548905829e-31"
